$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Product backlog text was re-ordered / split / re-worded (see commit msg:
# "Updated Backlog ... Backlog is now up to date with user stories").
# Rewrite B3:B24 with the new story text.
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "As a user, I want the alarm to vibrate,   `n    so I know when to switch activities."
$ws.Range("B4").Value = "As a user, I want the alarm to make a noise,`n`tSo I know when to switch activities.  "
$ws.Range("B5").Value = "As a user, I want to send feedback to the developers,`n    so they know how I feel about their app."
$ws.Range("B6").Value = "As a user, I want to report bugs I find, `n`tso they can be fixed quickly and efficiently."
$ws.Range("B7").Value = "As a user, I want to opt into and out of notifications to use the app,`n    so I can be reminded to use the application."
$ws.Range("B8").Value = "As a user, I want to be able to donate to the developers, `n    so I can reward them for creating a helpful application."
$ws.Range("B9").Value = "As a user, I want to choose different Pomodoro intervals,`nso I can manage my time how I want."
$ws.Range("B10").Value = "As a user, I want to be able to schedule my day,`n`tso I spend my time wisely."
$ws.Range("B11").Value = "As a user, I want to choose different alarms,`nso I can personilize it for me."
$ws.Range("B12").Value = "As a user, I want different color themes,`n`tso I can make it more appealing."
$ws.Range("B13").Value = "As a user, I want some motivational quotes,`n`tso I can motivate myself."
$ws.Range("B14").Value = "As a user, I want some motivational media, `n`tso I can motivate myself."
$ws.Range("B15").Value = "As a fitness coach, I want to share custom workouts,`n`tso that I can have my trainee follow it when I'm not there."
$ws.Range("B16").Value = "As a trainee, I want to see my profile's statistics,`n`tso that I can see and share how much I have improved over time."
$ws.Range("B17").Value = "As a user, I want to have a visible streak on my profile,`n`tso that I will be encouraged to not miss a day."
$ws.Range("B18").Value = "As a user, I want to favorite my preferred exercises on my profile,`n`tso that I can more easily find them later on."
$ws.Range("B19").Value = "As a workout buddy, I want to have the option to share on my profile when I've completed my scheduled workout,`n`tso that me and my partner can go through this experience together."
$ws.Range("B20").Value = "As a college student, I want to have access to exercises I can do in my chair,`n`tso that I can I can use this app in the library. "
$ws.Range("B21").Value = "As a desk-job-worker, I want stretches to be included as available exercises,`n`tso that I can keep my back, neck, and shoulders healthy during long days."
$ws.Range("B22").Value = "As a student, I want meditative exercises to be included in the app,`n`tso that I can manage my stress while remaining productive."
$ws.Range("B23").Value = "As a novice to fitness, I want visuals to help guide me during exercises,`n`tso that I can use the app without having to stop to research forms. "
$ws.Range("B24").Value = "As a novice to fitness, I want allow the app to randomly select exercises for me,`n`tso that I can see which ones work best for me. "

# Rows 22-24 previously had no user story (blank, centered style). They now
# hold wrapped story text like the rows above them, so copy the wrapped
# text format (style) from row 21 down onto them, then restore the ID values.
$ws.Range("B21").Copy() | Out-Null
$ws.Range("B22:B24").PasteSpecial(-4122) | Out-Null
$ws.Range("A22").Value = 21
$ws.Range("A23").Value = 22
$ws.Range("A24").Value = 23

# Those rows now wrap onto two lines like the others, so match their row
# height to the rest of the story rows.
$ws.Rows.Item(22).RowHeight = 31.5
$ws.Rows.Item(23).RowHeight = 31.5
$ws.Rows.Item(24).RowHeight = 31.5

# ---------------------------------------------------------------------------
# New "Sprint 1" rows: first three stories (rows 3, 4, 9) are now scheduled
# into Sprint 1, flagged "In Progress", worth 3 story points.
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = 1
$ws.Range("E3").Value = "In Progress"
$ws.Range("F3").Value = 3

$ws.Range("C4").Value = 1
$ws.Range("E4").Value = "In Progress"
$ws.Range("F4").Value = 3

$ws.Range("C9").Value = 1
$ws.Range("E9").Value = "In Progress"
$ws.Range("F9").Value = 3

# ---------------------------------------------------------------------------
# View state: selection moved to D3, and scroll reset to the top of the
# sheet (no frozen/forced top-left row anymore).
# ---------------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("D3").Select() | Out-Null
